# Natmi following Dr Hou advice
# Re-derive the Serpinf1-Plxdc1 LR-pair table to include the "ECs" sending
# cluster alongside the existing "FAPs" / "sCs" sending clusters, producing
# the full 3x3 sending/target cluster cross-join (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Serpinf1", "Plxdc1", "ECs", 2.0, 0.6666666666666666, 2.063913333333333, 6.19174, 0.003360194407478493, 0.003360194407478493, 3.0, 1.0, 2.932884333333333, 8.798653, 0.1407623943304514, 0.1407623943304514, 6.053219080691111, 54.47897172622, 0.0004729890102124651, 0.0004729890102124651),
    @("ECs", "Serpinf1", "Plxdc1", "FAPs", 2.0, 0.6666666666666666, 2.063913333333333, 6.19174, 0.003360194407478493, 0.003360194407478493, 3.0, 1.0, 4.225759666666667, 12.677279, 0.2028133335449359, 0.2028133335449359, 8.721601719495556, 78.49441547546, 0.0006814922291397637, 0.0006814922291397637),
    @("ECs", "Serpinf1", "Plxdc1", "sCs", 2.0, 0.6666666666666666, 2.063913333333333, 6.19174, 0.003360194407478493, 0.003360194407478493, 3.0, 1.0, 13.67706533333333, 41.031196, 0.6564242721246127, 0.6564242721246127, 28.22827750233778, 254.05449752104, 0.002205713168126264, 0.002205713168126264),
    @("FAPs", "Serpinf1", "Plxdc1", "ECs", 3.0, 1.0, 585.4210406666667, 1756.263122, 0.9531061576560254, 0.9531061576560252, 3.0, 1.0, 2.932884333333333, 8.798653, 0.1407623943304514, 0.1407623943304514, 1716.972198574963, 15452.74978717467, 0.1341615048027588, 0.1341615048027588),
    @("FAPs", "Serpinf1", "Plxdc1", "FAPs", 3.0, 1.0, 585.4210406666667, 1756.263122, 0.9531061576560254, 0.9531061576560252, 3.0, 1.0, 4.225759666666667, 12.677279, 0.2028133335449359, 0.2028133335449359, 2473.848621667227, 22264.63759500504, 0.1933026370564237, 0.1933026370564237),
    @("FAPs", "Serpinf1", "Plxdc1", "sCs", 3.0, 1.0, 585.4210406666667, 1756.263122, 0.9531061576560254, 0.9531061576560252, 3.0, 1.0, 13.67706533333333, 41.031196, 0.6564242721246127, 0.6564242721246127, 8006.841820705991, 72061.57638635392, 0.6256420157968429, 0.6256420157968428),
    @("sCs", "Serpinf1", "Plxdc1", "ECs", 3.0, 1.0, 26.739428, 80.218284, 0.04353364793649628, 0.04353364793649627, 3.0, 1.0, 2.932884333333333, 8.798653, 0.1407623943304514, 0.1407623943304514, 78.42364946349467, 705.812845171452, 0.006127900517480131, 0.00612790051748013),
    @("sCs", "Serpinf1", "Plxdc1", "FAPs", 3.0, 1.0, 26.739428, 80.218284, 0.04353364793649628, 0.04353364793649627, 3.0, 1.0, 4.225759666666667, 12.677279, 0.2028133335449359, 0.2028133335449359, 112.9943963521373, 1016.949567169236, 0.008829204259372428, 0.008829204259372427),
    @("sCs", "Serpinf1", "Plxdc1", "sCs", 3.0, 1.0, 26.739428, 80.218284, 0.04353364793649628, 0.04353364793649627, 3.0, 1.0, 13.67706533333333, 41.031196, 0.6564242721246127, 0.6564242721246127, 365.7169037319626, 3291.452133587664, 0.02857654315964372, 0.02857654315964371)
)

$numRows = $rows.Count
$numCols = $rows[0].Count
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

$ws.Range("A2:T10").Value = $data

Write-Host "Updated A2:T10 with the full sending/target cluster cross-join (9 rows)."
